$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 1
    3 = 2
    4 = 3
    5 = 8
    6 = 4
    7 = 3
    8 = 3
    9 = 3
    10 = 4
    11 = 9
    12 = 0
    13 = 5
    14 = 7
    15 = 3
    16 = 5
    17 = 7
    18 = 7
    19 = 5
    20 = 8
    21 = 3
    22 = 3
    23 = 5
    24 = 7
    25 = 4
    26 = 8
    27 = 2
    28 = 6
    29 = 6
    30 = 5
    31 = 5
    32 = 6
    33 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
